# se agrega ciclo para las lineas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values for rows 4 and 5 (MSISDN/CLIENTE/IMEI block),
# as if a loop re-wrote the line values with freshly generated numbers.
$ws.Range("A4").Value = "270670616"
$ws.Range("B4").Value = "3046010569"
$ws.Range("C4").Value = "883337485691834"

$ws.Range("A5").Value = "163908584"
$ws.Range("B5").Value = "3046010523"
$ws.Range("C5").Value = "883339511718342"

# Row 5 height matches the other data rows now (loop applies same height).
$ws.Range("A5:C5").RowHeight = 14.95

# C5 picks up the same formatting style that C4 already uses.
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C5").Value = "883339511718342"

# Move the active selection to C8, as left behind after the loop ran.
$ws.Range("C8").Select()
